# Generate Report for Handoff
# Mark the 665035ab-4a2e-4331-849a-c761baccb176 row as "Ready for handoff"
# on every sheet (Overview, zh-cn, de-de), and refresh the handoff
# timestamp recorded for that file on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B9").Value = "Ready for handoff"
$overview.Range("C9").Value = "Ready for handoff"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B9").Value = "Ready for handoff"
$zhcn.Range("D6").Value = "2016-03-09 02:30:34"
$zhcn.Range("D9").Value = "2016-03-09 02:30:34"
$zhcn.Range("D10").Value = "2016-03-09 02:30:34"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B9").Value = "Ready for handoff"
$dede.Range("D6").Value = "2016-03-09 02:30:45"
$dede.Range("D9").Value = "2016-03-09 02:30:45"
$dede.Range("D10").Value = "2016-03-09 02:30:45"
